# Update computed market/profit columns (H-N) across Leve profit sheets
# with refreshed values, as produced by the scheduled market-data runner.
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 881.25
$ws.Range("I12").Value = 284
$ws.Range("K12").Value = 284
$ws.Range("M12").Value = -114
$ws.Range("H17").Value = 2209.1638
$ws.Range("J17").Value = 2209.1638
$ws.Range("L17").Value = 6627.491399999999
$ws.Range("N17").Value = -6963.491399999999
$ws.Range("H28").Value = 651.1
$ws.Range("I28").Value = 177.73334
$ws.Range("J28").Value = 2071.2
$ws.Range("K28").Value = 177.73334
$ws.Range("L28").Value = 2071.2
$ws.Range("M28").Value = 307.26666
$ws.Range("N28").Value = -3041.2
$ws.Range("H107").Value = 1401.5714
$ws.Range("I107").Value = 1502.5
$ws.Range("J107").Value = 1361.2
$ws.Range("K107").Value = 1502.5
$ws.Range("L107").Value = 1361.2
$ws.Range("M107").Value = 417.5
$ws.Range("N107").Value = -5201.2
$ws.Range("H113").Value = 2949.1765
$ws.Range("I113").Value = 2413.75
$ws.Range("J113").Value = 3113.923
$ws.Range("K113").Value = 2413.75
$ws.Range("L113").Value = 3113.923
$ws.Range("M113").Value = 840.25
$ws.Range("N113").Value = -9621.922999999999
$ws.Range("H130").Value = 49348.215
$ws.Range("J130").Value = 49348.215
$ws.Range("L130").Value = 49348.215
$ws.Range("N130").Value = -59388.215
$ws.Range("H132").Value = 1662.1111
$ws.Range("I132").Value = 1750.225
$ws.Range("J132").Value = 957.2
$ws.Range("K132").Value = 5250.674999999999
$ws.Range("L132").Value = 2871.6
$ws.Range("M132").Value = -2720.674999999999
$ws.Range("N132").Value = -7931.6
$ws.Range("H137").Value = 2252.262
$ws.Range("I137").Value = 1796.0714
$ws.Range("J137").Value = 3164.6428
$ws.Range("K137").Value = 5388.2142
$ws.Range("L137").Value = 9493.928400000001
$ws.Range("M137").Value = -2838.2142
$ws.Range("N137").Value = -14593.9284

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1766.5807
$ws.Range("I45").Value = 1731.0435
$ws.Range("J45").Value = 1868.75
$ws.Range("K45").Value = 1731.0435
$ws.Range("L45").Value = 1868.75
$ws.Range("M45").Value = -1354.0435
$ws.Range("N45").Value = -2622.75
$ws.Range("H48").Value = 0
$ws.Range("J48").Value = 0
$ws.Range("L48").Value = 0
$ws.Range("N48").Value = $null
$ws.Range("H61").Value = 6634.9365
$ws.Range("I61").Value = 3924.32
$ws.Range("J61").Value = 17060.385
$ws.Range("K61").Value = 3924.32
$ws.Range("L61").Value = 17060.385
$ws.Range("M61").Value = -3712.32
$ws.Range("N61").Value = -17484.385
$ws.Range("H74").Value = 6662.683
$ws.Range("I74").Value = 4588.5
$ws.Range("K74").Value = 4588.5
$ws.Range("M74").Value = -3714.5
$ws.Range("H77").Value = 6662.683
$ws.Range("I77").Value = 4588.5
$ws.Range("K77").Value = 22942.5
$ws.Range("M77").Value = -18574.5
$ws.Range("H110").Value = 0
$ws.Range("I110").Value = 0
$ws.Range("K110").Value = 0
$ws.Range("M110").Value = $null
$ws.Range("H132").Value = 4655.905
$ws.Range("I132").Value = 1646.5217
$ws.Range("J132").Value = 8298.842000000001
$ws.Range("K132").Value = 4939.5651
$ws.Range("L132").Value = 24896.526
$ws.Range("M132").Value = -2409.5651
$ws.Range("N132").Value = -29956.526
$ws.Range("H136").Value = 6634.9365
$ws.Range("I136").Value = 3924.32
$ws.Range("J136").Value = 17060.385
$ws.Range("K136").Value = 11772.96
$ws.Range("L136").Value = 51181.155
$ws.Range("M136").Value = -9222.960000000001
$ws.Range("N136").Value = -56281.155

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H133").Value = 52415
$ws.Range("J133").Value = 52415
$ws.Range("L133").Value = 52415
$ws.Range("N133").Value = -62535
$ws.Range("H134").Value = 20233.309
$ws.Range("I134").Value = 1960.5897
$ws.Range("J134").Value = 64773.062
$ws.Range("K134").Value = 5881.7691
$ws.Range("L134").Value = 194319.186
$ws.Range("M134").Value = -3346.7691
$ws.Range("N134").Value = -199389.186

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2151.2322
$ws.Range("I31").Value = 1639.9487
$ws.Range("J31").Value = 3324.1765
$ws.Range("K31").Value = 1639.9487
$ws.Range("L31").Value = 3324.1765
$ws.Range("M31").Value = -1344.9487
$ws.Range("N31").Value = -3914.1765
$ws.Range("H34").Value = 2151.2322
$ws.Range("I34").Value = 1639.9487
$ws.Range("J34").Value = 3324.1765
$ws.Range("K34").Value = 1639.9487
$ws.Range("L34").Value = 3324.1765
$ws.Range("M34").Value = -1437.9487
$ws.Range("N34").Value = -3728.1765
$ws.Range("H94").Value = 1342.3334
$ws.Range("I94").Value = 1012
$ws.Range("K94").Value = 1012
$ws.Range("M94").Value = -561
$ws.Range("H99").Value = 2969.375
$ws.Range("I99").Value = 2400
$ws.Range("J99").Value = 3311
$ws.Range("K99").Value = 2400
$ws.Range("L99").Value = 3311
$ws.Range("M99").Value = -902
$ws.Range("N99").Value = -6307
$ws.Range("H126").Value = 2969.375
$ws.Range("I126").Value = 2400
$ws.Range("J126").Value = 3311
$ws.Range("K126").Value = 7200
$ws.Range("L126").Value = 9933
$ws.Range("M126").Value = -4730
$ws.Range("N126").Value = -14873
$ws.Range("H134").Value = 2957.2354
$ws.Range("I134").Value = 1947.24
$ws.Range("J134").Value = 3928.3845
$ws.Range("K134").Value = 5841.72
$ws.Range("L134").Value = 11785.1535
$ws.Range("M134").Value = -3306.72
$ws.Range("N134").Value = -16855.1535

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 7937289
$ws.Range("I5").Value = 443.875
$ws.Range("J5").Value = 33335194
$ws.Range("K5").Value = 1331.625
$ws.Range("L5").Value = 100005582
$ws.Range("M5").Value = -1219.625
$ws.Range("N5").Value = -100005806
$ws.Range("H80").Value = 2484.6155
$ws.Range("I80").Value = 2666.6667
$ws.Range("J80").Value = 2430
$ws.Range("K80").Value = 8000.000100000001
$ws.Range("L80").Value = 7290
$ws.Range("M80").Value = -7064.000100000001
$ws.Range("N80").Value = -9162
$ws.Range("H83").Value = 2484.6155
$ws.Range("I83").Value = 2666.6667
$ws.Range("J83").Value = 2430
$ws.Range("K83").Value = 24000.0003
$ws.Range("L83").Value = 21870
$ws.Range("M83").Value = -19320.0003
$ws.Range("N83").Value = -31230
$ws.Range("H107").Value = 1417.7858
$ws.Range("I107").Value = 324.46155
$ws.Range("J107").Value = 2365.3333
$ws.Range("K107").Value = 973.38465
$ws.Range("L107").Value = 7095.999899999999
$ws.Range("M107").Value = 946.61535
$ws.Range("N107").Value = -10935.9999
$ws.Range("H131").Value = 41071.793
$ws.Range("I131").Value = 2150
$ws.Range("J131").Value = 60532.688
$ws.Range("K131").Value = 6450
$ws.Range("L131").Value = 181598.064
$ws.Range("M131").Value = -1410
$ws.Range("N131").Value = -191678.064
$ws.Range("H135").Value = 7937289
$ws.Range("I135").Value = 443.875
$ws.Range("J135").Value = 33335194
$ws.Range("K135").Value = 3994.875
$ws.Range("L135").Value = 300016746
$ws.Range("M135").Value = -1459.875
$ws.Range("N135").Value = -300021816

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1782.8334
$ws.Range("I113").Value = 1655.5294
$ws.Range("J113").Value = 1896.7368
$ws.Range("K113").Value = 1655.5294
$ws.Range("L113").Value = 1896.7368
$ws.Range("M113").Value = 514.4706000000001
$ws.Range("N113").Value = -6236.7368

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3923.5264
$ws.Range("I40").Value = 3766.077
$ws.Range("K40").Value = 3766.077
$ws.Range("M40").Value = -3630.077
$ws.Range("H76").Value = 27716
$ws.Range("J76").Value = 27716
$ws.Range("L76").Value = 27716
$ws.Range("N76").Value = -28392
$ws.Range("H79").Value = 27716
$ws.Range("J79").Value = 27716
$ws.Range("L79").Value = 27716
$ws.Range("N79").Value = -30056
$ws.Range("H122").Value = 6952.9
$ws.Range("I122").Value = 6565.8096
$ws.Range("J122").Value = 7856.1113
$ws.Range("K122").Value = 19697.4288
$ws.Range("L122").Value = 23568.3339
$ws.Range("M122").Value = -17247.4288
$ws.Range("N122").Value = -28468.3339
$ws.Range("H132").Value = 3367.4626
$ws.Range("I132").Value = 3453.843
$ws.Range("J132").Value = 3092.125
$ws.Range("K132").Value = 10361.529
$ws.Range("L132").Value = 9276.375
$ws.Range("M132").Value = -7831.528999999999
$ws.Range("N132").Value = -14336.375
$ws.Range("H135").Value = 41766.668
$ws.Range("J135").Value = 41766.668
$ws.Range("L135").Value = 41766.668
$ws.Range("N135").Value = -51906.668
$ws.Range("H136").Value = 3708.4153
$ws.Range("I136").Value = 2125.6956
$ws.Range("K136").Value = 6377.0868
$ws.Range("M136").Value = -3827.0868

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("L95").Value = 0
$ws.Range("N95").Value = $null
$ws.Range("H126").Value = 908.1875
$ws.Range("I126").Value = 791.5
$ws.Range("J126").Value = 1725
$ws.Range("K126").Value = 2374.5
$ws.Range("L126").Value = 5175
$ws.Range("M126").Value = 95.5
$ws.Range("N126").Value = -10115
$ws.Range("H136").Value = 3086.9524
$ws.Range("I136").Value = 1340.6608
$ws.Range("J136").Value = 6579.5356
$ws.Range("K136").Value = 4021.9824
$ws.Range("L136").Value = 19738.6068
$ws.Range("M136").Value = -1471.9824
$ws.Range("N136").Value = -24838.6068
